# Chiffres COVID-19 Valais - apply the data corrections for rows 421-427
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Row 421: "Nb nouveaux cas positifs" corrected 104 -> 103 ---
$ws.Range("C421").Value = 103

# --- Row 423: "Nb nouveaux cas positifs" corrected 108 -> 107 ---
$ws.Range("C423").Value = 107

# --- Row 425: "Nb nouveaux cas positifs" corrected 48 -> 52 ---
$ws.Range("C425").Value = 52

# --- Row 426: "Nb nouveaux cas positifs" corrected 9 -> 76 ---
$ws.Range("C426").Value = 76

# --- Row 427: newly filled-in day (was blank/placeholder before) ---
$ws.Range("C427").Value = 10
$ws.Range("E427").Value = 8
$ws.Range("F427").Value = 7
$ws.Range("G427").Value = 29

# L427 and M427 are formatted as Text (number format "@"). Assigning a plain
# numeric .Value to such a cell stores it as a shared-string instead of a
# real number, and changing .NumberFormat directly allocates brand-new style
# records. To land a genuine numeric 0 while keeping the existing style
# (and without growing the style table), temporarily borrow the number
# format from a General-formatted cell with the same border via
# PasteSpecial (formats only), set the numeric value, then restore the
# original Text format the same way.
$ws.Range("C3").Copy() | Out-Null
$ws.Range("L427").PasteSpecial(-4122) | Out-Null
$ws.Range("L427").Value = 0
$ws.Range("L421").Copy() | Out-Null
$ws.Range("L427").PasteSpecial(-4122) | Out-Null

$ws.Range("C144").Copy() | Out-Null
$ws.Range("M427").PasteSpecial(-4122) | Out-Null
$ws.Range("M427").Value = 0
$ws.Range("M421").Copy() | Out-Null
$ws.Range("M427").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false
